$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the small "source folder" cross-reference table in G4:H5
# (I2C -> JA, I2S -> JB), reusing the body-row format (font/fill) from
# an existing labeled cell but with a left-only thin border.
$ws.Range("E4").Copy()
foreach ($addr in "G4", "H4", "G5", "H5") {
    $cell = $ws.Range($addr)
    $cell.PasteSpecial(-4122)
    $cell.Borders.LineStyle = -4142
    $cell.Borders.Item(7).LineStyle = 1
}

$ws.Range("G4").Value = "I2C"
$ws.Range("H4").Value = "JA"
$ws.Range("G5").Value = "I2S"
$ws.Range("H5").Value = "JB"

# Move the selection (matches the new cursor position left after editing)
$ws.Range("J14").Select()
